$p = $ppt.ActivePresentation

# Slide 12 ("logboek" week overview) contains the schedule table "Table 8".
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item("Table 8")
$tbl = $shp.Table

# Row 7 ("Tjorben Godeau"), column 5 ("17u00-18u00") currently reads
# "Clash of Clans" - update it to "Testing".
$cell = $tbl.Cell(7, 5)
$cell.Shape.TextFrame.TextRange.Text = "Testing"
